$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Replace the greeting text in E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Leave E8 selected, as it was the active cell when the file was last saved
$ws.Range("E8").Select()
